$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.118.40"
$ws.Range("E2").Value = "  +0.57%  "

# Row 3
$ws.Range("D3").Value = "'1.894.73"
$ws.Range("E3").Value = "  +0.45%  "

# Row 4
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.27%  "

# Row 5
$ws.Range("D5").Value = "'323.72"
$ws.Range("E5").Value = "  -1.78%  "

# Row 6
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.25%  "

# Row 7
$ws.Range("D7").Value = "'0.4698"
$ws.Range("E7").Value = "  +2.48%  "

# Row 8
$ws.Range("E8").Value = "  -2.36%  "

# Row 9
$ws.Range("D9").Value = "'47.41"
$ws.Range("E9").Value = "  -0.69%  "

# Row 10
$ws.Range("D10").Value = "'0.07982"
$ws.Range("E10").Value = "  +0.36%  "

# Row 11
$ws.Range("D11").Value = "'0.9900"
$ws.Range("E11").Value = "  -0.53%  "

# Row 12
$ws.Range("D12").Value = "'22.45"
$ws.Range("E12").Value = "  +3.61%  "

# Row 13
$ws.Range("D13").Value = "'1.929.84"
$ws.Range("E13").Value = "  +2.42%  "

# Row 14
$ws.Range("D14").Value = "'5.840"
$ws.Range("E14").Value = "  -1.27%  "

# Row 15
$ws.Range("D15").Value = "'7.022"
$ws.Range("E15").Value = "  -0.64%  "

# Row 16
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.28%  "

# Row 17
$ws.Range("D17").Value = "'88.52"
$ws.Range("E17").Value = "  -0.01%  "

# Row 18
$ws.Range("D18").Value = "'0.06615"
$ws.Range("E18").Value = "  +0.75%  "

# Row 19
$ws.Range("E19").Value = "  -0.19%  "

# Row 20
$ws.Range("D20").Value = "'17.42"
$ws.Range("E20").Value = "  +0.15%  "

# Row 21
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  +0.20%  "

# Row 22
$ws.Range("D22").Value = "'29.140.77"
$ws.Range("E22").Value = "  +0.63%  "

# Row 23
$ws.Range("D23").Value = "'5.486"
$ws.Range("E23").Value = "  +1.13%  "

# Row 24
$ws.Range("D24").Value = "'11.47"
$ws.Range("E24").Value = "  +0.29%  "

# Row 25
$ws.Range("D25").Value = "'2.200"
$ws.Range("E25").Value = "  +0.15%  "

# Row 26
$ws.Range("D26").Value = "'2.126.83"
$ws.Range("E26").Value = "  +0.80%  "

# Row 27
$ws.Range("D27").Value = "'153.92"
$ws.Range("E27").Value = "  -1.60%  "

# Row 28
$ws.Range("E28").Value = "  +0.18%  "

# Row 29
$ws.Range("D29").Value = "'5.995"
$ws.Range("E29").Value = "  +9.33%  "

# Row 30
$ws.Range("D30").Value = "'2.077"
$ws.Range("E30").Value = "  -0.29%  "

# Row 31
$ws.Range("D31").Value = "'116.83"
$ws.Range("E31").Value = "  -0.57%  "

# Row 32
$ws.Range("D32").Value = "'1.054"
$ws.Range("E32").Value = "  +0.97%  "

# Row 33
$ws.Range("D33").Value = "'0.09434"
$ws.Range("E33").Value = "  +1.26%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.553"
$ws.Range("E34").Value = "  +0.74%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.393"
$ws.Range("E35").Value = "  -0.72%  "

# Row 36
$ws.Range("D36").Value = "'5.333"

# Row 37
$ws.Range("D37").Value = "'0.06056"
$ws.Range("E37").Value = "  -0.11%  "

# Row 38
$ws.Range("E38").Value = "  +0.51%  "

# Row 39
$ws.Range("D39").Value = "'1.170"
$ws.Range("E39").Value = "  -0.07%  "

# Row 40
$ws.Range("D40").Value = "'8.047"
$ws.Range("E40").Value = "  -3.44%  "

# Row 41
$ws.Range("D41").Value = "'0.5791"
$ws.Range("E41").Value = "  +0.10%  "

# Row 42
$ws.Range("D42").Value = "'0.1819"
$ws.Range("E42").Value = "  -0.15%  "

# Row 43
$ws.Range("D43").Value = "'2.450"
$ws.Range("E43").Value = "  +7.95%  "

# Row 44
$ws.Range("D44").Value = "'10.03"
$ws.Range("E44").Value = "  -0.48%  "

# Row 45
$ws.Range("D45").Value = "'0.07675"
$ws.Range("E45").Value = "  +2.02%  "

# Row 46
$ws.Range("D46").Value = "'1.251"
$ws.Range("E46").Value = "  -0.60%  "

# Row 47
$ws.Range("D47").Value = "'12.08"
$ws.Range("E47").Value = "  +0.49%  "

# Row 48
$ws.Range("D48").Value = "'0.5449"
$ws.Range("E48").Value = "  -0.04%  "

# Row 49
$ws.Range("D49").Value = "'1.895"
$ws.Range("E49").Value = "  -0.29%  "

# Row 50
$ws.Range("D50").Value = "'113.17"
$ws.Range("E50").Value = "  +1.72%  "

# Row 51
$ws.Range("D51").Value = "'43.58"
$ws.Range("E51").Value = "  -2.26%  "
